$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("W2").Value = 2.34
$ws.Range("J3").Value = 3.1
$ws.Range("S3").Value = 3.95
$ws.Range("F4").Value = 1.8
$ws.Range("G4").Value = 1.93
$ws.Range("H4").Value = 4.4
$ws.Range("V4").Value = 1.21
$ws.Range("W4").Value = 2.06
$ws.Range("I5").Value = 1.55
$ws.Range("Q5").Value = 1.79
$ws.Range("R5").Value = 1.47
$ws.Range("V5").Value = 2.8
$ws.Range("Z5").Value = 9.199999999999999
$ws.Range("AA5").Value = 13.5
$ws.Range("AB5").Value = 25
$ws.Range("AI5").Value = 34
$ws.Range("H6").Value = 12
$ws.Range("I6").Value = 13
$ws.Range("J6").Value = 7.8
$ws.Range("K6").Value = 8
$ws.Range("P6").Value = 3.3
$ws.Range("Q6").Value = 1.4
$ws.Range("R6").Value = 1.93
$ws.Range("T6").Value = 1.86
$ws.Range("Y6").Value = 55
$ws.Range("AA6").Value = 580
$ws.Range("AB6").Value = 13.5
$ws.Range("AD6").Value = 44
$ws.Range("AH6").Value = 28
$ws.Range("AI6").Value = 120
$ws.Range("AJ6").Value = 10.5
$ws.Range("AL6").Value = 28
$ws.Range("AN6").Value = 3.4
$ws.Range("AO6").Value = 140
$ws.Range("F7").Value = 3.5
$ws.Range("G7").Value = 3.6
$ws.Range("H7").Value = 2.2
$ws.Range("I7").Value = 2.24
$ws.Range("T7").Value = 1.63
$ws.Range("U7").Value = 2.48
$ws.Range("V7").Value = 1.81
$ws.Range("W7").Value = 1.38
$ws.Range("X7").Value = 17.5
$ws.Range("Y7").Value = 12.5
$ws.Range("AA7").Value = 28
$ws.Range("AJ7").Value = 60
$ws.Range("AK7").Value = 36
$ws.Range("AL7").Value = 40
$ws.Range("AN7").Value = 28
$ws.Range("F8").Value = 1.47
$ws.Range("G8").Value = 1.49
$ws.Range("I8").Value = 8.199999999999999
$ws.Range("K8").Value = 5.3
$ws.Range("N8").Value = 4.4
$ws.Range("P8").Value = 2.14
$ws.Range("U8").Value = 1.88
$ws.Range("V8").Value = 1.14
$ws.Range("Y8").Value = 26
$ws.Range("Z8").Value = 70
$ws.Range("AA8").Value = 280
$ws.Range("AC8").Value = 11
$ws.Range("AE8").Value = 130
$ws.Range("AI8").Value = 120
$ws.Range("AM8").Value = 160
$ws.Range("AO8").Value = 160
$ws.Range("F9").Value = 3.25
$ws.Range("G9").Value = 3.35
$ws.Range("J9").Value = 4
$ws.Range("K9").Value = 4.1
$ws.Range("M9").Value = 1.04
$ws.Range("P9").Value = 2.58
$ws.Range("Q9").Value = 1.59
$ws.Range("R9").Value = 1.64
$ws.Range("S9").Value = 2.48
$ws.Range("W9").Value = 1.42
$ws.Range("Y9").Value = 14.5
$ws.Range("Z9").Value = 17
$ws.Range("AK9").Value = 32
$ws.Range("AM9").Value = 60
$ws.Range("AN9").Value = 21
$ws.Range("G10").Value = 2.34
$ws.Range("H10").Value = 3.1
$ws.Range("U10").Value = 2.78
$ws.Range("W10").Value = 1.75
$ws.Range("X10").Value = 25
$ws.Range("Z10").Value = 26
$ws.Range("AB10").Value = 16
$ws.Range("AC10").Value = 9.4
$ws.Range("AD10").Value = 13.5
$ws.Range("AE10").Value = 29
$ws.Range("AG10").Value = 11.5
$ws.Range("AH10").Value = 14
$ws.Range("AI10").Value = 32
$ws.Range("AK10").Value = 20
$ws.Range("AL10").Value = 26
$ws.Range("AN10").Value = 11
$ws.Range("AO10").Value = 18
$ws.Range("F11").Value = 2.24
$ws.Range("H11").Value = 3.35
$ws.Range("J11").Value = 3.75
$ws.Range("Q11").Value = 1.69
$ws.Range("S11").Value = 2.78
$ws.Range("T11").Value = 1.63
$ws.Range("U11").Value = 2.48
$ws.Range("AB11").Value = 12.5
$ws.Range("AC11").Value = 8.800000000000001
$ws.Range("AH11").Value = 15.5
$ws.Range("AJ11").Value = 28
$ws.Range("AM11").Value = 70
$ws.Range("AN11").Value = 13
$ws.Range("I12").Value = 18.5
$ws.Range("J12").Value = 10
$ws.Range("P12").Value = 4.3
$ws.Range("Q12").Value = 1.28
$ws.Range("T12").Value = 1.8
$ws.Range("U12").Value = 2.16
$ws.Range("AB12").Value = 19
$ws.Range("AG12").Value = 13.5
$ws.Range("AL12").Value = 30
$ws.Range("AN12").Value = 2.48
$ws.Range("F13").Value = 3.15
$ws.Range("H13").Value = 2.34
$ws.Range("I13").Value = 2.38
$ws.Range("L13").Value = 1.32
$ws.Range("N13").Value = 5
$ws.Range("P13").Value = 2.32
$ws.Range("Q13").Value = 1.71
$ws.Range("R13").Value = 1.54
$ws.Range("S13").Value = 2.76
$ws.Range("T13").Value = 1.58
$ws.Range("U13").Value = 2.5
$ws.Range("V13").Value = 1.72
$ws.Range("W13").Value = 1.45
$ws.Range("X13").Value = 18.5
$ws.Range("Y13").Value = 13.5
$ws.Range("Z13").Value = 17
$ws.Range("AB13").Value = 16
$ws.Range("AH13").Value = 15
$ws.Range("AI13").Value = 30
$ws.Range("AJ13").Value = 55
$ws.Range("AK13").Value = 32
$ws.Range("AL13").Value = 38
$ws.Range("AM13").Value = 65
$ws.Range("AN13").Value = 23
$ws.Range("AO13").Value = 17
